# Append/update timestamps on the "ランサーズ" sheet:
# rows 2-8, column A ("取得日時") change from 2025-10-24 06:27:19
# to 2025-10-24 06:33:35 (new scrape run timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-24 06:33:35"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
